$p = $ppt.ActivePresentation

# Slide 3, "Text Placeholder 17" shape holds the bullet:
#   "California has second highest sales at $15.39 Million"
# The figure "$15.39 Million" needs to be corrected to "$15.30 Million".
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(3)
$tf = $shp.TextFrame
$tr = $tf.TextRange

$fullText = $tr.Text
$target = "at `$15.39 Million"
$startIdx = $fullText.IndexOf($target) + 1   # TextRange.Characters is 1-based

if ($startIdx -gt 0) {
    # Re-type just the dollar figure ("$15.39 " -> "$15.30 "), leaving "at " and
    # "Million" as their own runs (PowerPoint splits the original run at the
    # edited boundary when only part of it is replaced).
    $numStart = $startIdx + 3   # skip past "at "
    $numRange = $tr.Characters($numStart, 7)   # "$15.39 " / "$15.30 " are 7 chars
    $numRange.Text = "`$15.30 "
}
